# Update the build version / timestamp string throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldStamp)"
$newVersion = "mines - January 30 (built on $newStamp)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Ashton Coal Mine, Australia, M0007, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds "build_version"; header is row 1, data rows 2..17
$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 17 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}

$wb.Save()
